$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.727.63"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.454.84"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "576.11"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "160.86"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.16%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +14.57%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "3.458.07"
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  -1.63%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.125"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.449"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.95%  "
$ws.Range("D13").Value = "4.052.57"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("E15").Value = "  +0.27%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "28.20"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "64.799.99"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "3.473.90"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +3.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "380.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "8.11"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.552"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +4.34%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "72.68"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("E26").Value = "  +0.64%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.97"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +6.60%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +11.48%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.11"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("E32").Value = "  +2.01%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "23.62"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("E34").Value = "  +5.99%  "
$ws.Range("E35").Value = "  +11.57%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "161.28"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").Value = "2.953.04"
$ws.Range("E39").Value = "  -0.88%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "26.60"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.00%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.60"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("E43").Value = "  +1.95%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "42.73"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.774"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.17%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "25.65"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +11.18%  "
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("E48").Value = "  +8.62%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "308.56"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +5.57%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.868"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +4.65%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.62"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.78%  "
